# added more games, sped up simulate game logic, and drafted optimization logic
# -> refreshed the Longwood_A team transition-probability matrix (Sheet1)
#    with updated simulation results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2202970297029703
$ws.Range("C2").Value2 = 0.4826732673267327
$ws.Range("J2").Value2 = 0.02227722772277228
$ws.Range("P2").Value2 = 0.1584158415841584
$ws.Range("S2").Value2 = 0.1163366336633663
$ws.Range("B3").Value2 = 0.01
$ws.Range("C3").Value2 = 0.005
$ws.Range("J3").Value2 = 0.065
$ws.Range("O3").Value2 = 0.005
$ws.Range("P3").Value2 = 0.735
$ws.Range("S3").Value2 = 0.18
$ws.Range("J4").Value2 = 0.08
$ws.Range("P4").Value2 = 0.66
$ws.Range("S4").Value2 = 0.26
$ws.Range("B6").Value2 = 0.05737704918032787
$ws.Range("D6").Value2 = 0.01229508196721311
$ws.Range("F6").Value2 = 0.09016393442622951
$ws.Range("J6").Value2 = 0.3319672131147541
$ws.Range("O6").Value2 = 0.01639344262295082
$ws.Range("Q6").Value2 = 0.1188524590163934
$ws.Range("R6").Value2 = 0.06967213114754098
$ws.Range("S6").Value2 = 0.3032786885245902
$ws.Range("B7").Value2 = 0.1573033707865168
$ws.Range("D7").Value2 = 0.01123595505617977
$ws.Range("F7").Value2 = 0.06179775280898876
$ws.Range("J7").Value2 = 0.1460674157303371
$ws.Range("O7").Value2 = 0.02808988764044944
$ws.Range("Q7").Value2 = 0.1797752808988764
$ws.Range("R7").Value2 = 0.06179775280898876
$ws.Range("S7").Value2 = 0.3539325842696629
$ws.Range("B8").Value2 = 0.1201923076923077
$ws.Range("D8").Value2 = 0.02163461538461538
$ws.Range("F8").Value2 = 0.0625
$ws.Range("J8").Value2 = 0.1370192307692308
$ws.Range("O8").Value2 = 0.02644230769230769
$ws.Range("Q8").Value2 = 0.1923076923076923
$ws.Range("R8").Value2 = 0.07692307692307693
$ws.Range("S8").Value2 = 0.3629807692307692
$ws.Range("B9").Value2 = 0.1027027027027027
$ws.Range("D9").Value2 = 0.01621621621621622
$ws.Range("F9").Value2 = 0.04864864864864865
$ws.Range("J9").Value2 = 0.1513513513513514
$ws.Range("O9").Value2 = 0.02702702702702703
$ws.Range("Q9").Value2 = 0.2054054054054054
$ws.Range("R9").Value2 = 0.06486486486486487
$ws.Range("S9").Value2 = 0.3837837837837838
$ws.Range("B10").Value2 = 0.1325549450549451
$ws.Range("D10").Value2 = 0.02541208791208791
$ws.Range("F10").Value2 = 0.06043956043956044
$ws.Range("J10").Value2 = 0.1476648351648352
$ws.Range("O10").Value2 = 0.01923076923076923
$ws.Range("Q10").Value2 = 0.2438186813186813
$ws.Range("R10").Value2 = 0.0570054945054945
$ws.Range("S10").Value2 = 0.3138736263736264
$ws.Range("G11").Value2 = 0.1486486486486487
$ws.Range("J11").Value2 = 0.1148648648648649
$ws.Range("K11").Value2 = 0.1891891891891892
$ws.Range("L11").Value2 = 0.5236486486486487
$ws.Range("S11").Value2 = 0.02364864864864865
$ws.Range("F12").Value2 = 0.006097560975609756
$ws.Range("G12").Value2 = 0.6707317073170732
$ws.Range("J12").Value2 = 0.274390243902439
$ws.Range("K12").Value2 = 0.006097560975609756
$ws.Range("L12").Value2 = 0.02439024390243903
$ws.Range("S12").Value2 = 0.01829268292682927
$ws.Range("G13").Value2 = 0.5833333333333334
$ws.Range("J13").Value2 = 0.375
$ws.Range("S13").Value2 = 0.04166666666666666
$ws.Range("F15").Value2 = 0.02158273381294964
$ws.Range("H15").Value2 = 0.1690647482014389
$ws.Range("I15").Value2 = 0.0683453237410072
$ws.Range("J15").Value2 = 0.3633093525179856
$ws.Range("K15").Value2 = 0.07913669064748201
$ws.Range("M15").Value2 = 0.007194244604316547
$ws.Range("N15").Value2 = 0.003597122302158274
$ws.Range("O15").Value2 = 0.07913669064748201
$ws.Range("S15").Value2 = 0.2086330935251799
$ws.Range("F16").Value2 = 0.02531645569620253
$ws.Range("H16").Value2 = 0.1687763713080169
$ws.Range("I16").Value2 = 0.1139240506329114
$ws.Range("J16").Value2 = 0.4008438818565401
$ws.Range("K16").Value2 = 0.1181434599156118
$ws.Range("M16").Value2 = 0.01265822784810127
$ws.Range("N16").Value2 = 0.004219409282700422
$ws.Range("O16").Value2 = 0.06329113924050633
$ws.Range("S16").Value2 = 0.09282700421940929
$ws.Range("F17").Value2 = 0.03076923076923077
$ws.Range("H17").Value2 = 0.1653846153846154
$ws.Range("I17").Value2 = 0.075
$ws.Range("J17").Value2 = 0.4480769230769231
$ws.Range("K17").Value2 = 0.07692307692307693
$ws.Range("M17").Value2 = 0.01923076923076923
$ws.Range("O17").Value2 = 0.08461538461538462
$ws.Range("S17").Value2 = 0.1
$ws.Range("F18").Value2 = 0.01973684210526316
$ws.Range("H18").Value2 = 0.1776315789473684
$ws.Range("I18").Value2 = 0.07236842105263158
$ws.Range("J18").Value2 = 0.4342105263157895
$ws.Range("K18").Value2 = 0.09210526315789473
$ws.Range("M18").Value2 = 0.0131578947368421
$ws.Range("O18").Value2 = 0.125
$ws.Range("S18").Value2 = 0.06578947368421052
$ws.Range("F19").Value2 = 0.02327586206896552
$ws.Range("H19").Value2 = 0.1913793103448276
$ws.Range("I19").Value2 = 0.075
$ws.Range("J19").Value2 = 0.3862068965517241
$ws.Range("K19").Value2 = 0.1103448275862069
$ws.Range("M19").Value2 = 0.03103448275862069
$ws.Range("O19").Value2 = 0.08879310344827586
$ws.Range("S19").Value2 = 0.0939655172413793
